# Major revision of interim report: fix "discounting" sheet formulas so the
# NPV-style calculations (AO:AX) discount the capital cost base (AM, "total
# fixed+variable capital") rather than the stray AN column, recompute the
# dependent totals/ratios, and update the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("discounting")
$ws.Activate()

# --- Fix the discounting formulas: reference $AM (capital cost) instead of $AN ---

# Row 3 block (AO3 is a standalone formula, AP3:AX3 is a shared block)
$ws.Range("AO3").Formula = "=`$AM3/10*1/(POWER(1.08,AO`$2-1))"
$ws.Range("AP3:AX3").Formula = "=`$AM3/10*1/(POWER(1.08,AP`$2-1))"

# Rows 4-16 share one formula block anchored at AO4
$ws.Range("AO4:AX16").Formula = "=`$AM4/10*1/(POWER(1.08,AO`$2-1))"

# --- Recompute / re-enter the dependent total column (AY) ---
$ws.Range("AY4").Formula = "=SUM(AO4:AX4)"
$ws.Range("AY5:AY16").Formula = "=SUM(AO5:AX5)"

# --- A few other cells were re-entered (same result, broken out of their
#     previous shared-formula group) ---
$ws.Range("Q6").Formula = "=SUM(G6:P6)"
$ws.Range("Q11").Formula = "=SUM(G11:P11)"
$ws.Range("BB7").Formula = "=AH7-AY7"

# --- Update the saved selection on the sheet ---
$ws.Range("K2").Select()
